# Regenerate the Overview / zh-cn / de-de handoff status rows.
# Row order becomes [ffff96a6b1e0, ffffffdf71933c, 9ac5998c] and the
# 9ac5998c entry flips from "Handed back" to "Ready for handoff" with
# fresh timestamps, matching the regenerated CI report.
$wb = $excel.ActiveWorkbook

# ---- Sheet: Overview ----
$ws = $wb.Worksheets.Item("Overview")

# Update cell values (data rows 2-4; header row unchanged)
$ws.Range("A2").Value = "ffff96a6b1e0-7098-44cb-bf6b-eed6400b361f.md"
$ws.Range("B2").Value = "Handed back: in sync with en-US"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "2016-03-24 03:14:40"
$ws.Range("A3").Value = "ffffffdf71933c-d945-47f1-b0b7-c48448e4a22c.md"
$ws.Range("B3").Value = "Handed back: in sync with en-US"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "2016-03-24 03:14:40"
$ws.Range("A4").Value = "9ac5998c-9398-47e2-b777-121c530ec423.md"
$ws.Range("B4").Value = "Ready for handoff"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "2016-03-24 03:19:26"

# Update hyperlink display text in place (keep existing target URL/rId)
foreach ($h in $ws.Hyperlinks) {
  $addr = $h.Range.Address()
  if ($addr -eq '$A$2') { $h.TextToDisplay = "ffff96a6b1e0-7098-44cb-bf6b-eed6400b361f.md" }
  if ($addr -eq '$A$3') { $h.TextToDisplay = "ffffffdf71933c-d945-47f1-b0b7-c48448e4a22c.md" }
  if ($addr -eq '$A$4') { $h.TextToDisplay = "9ac5998c-9398-47e2-b777-121c530ec423.md" }
}

# ---- Sheet: zh-cn ----
$ws = $wb.Worksheets.Item("zh-cn")

# Update cell values (data rows 2-4; header row unchanged)
$ws.Range("A2").Value = "ffff96a6b1e0-7098-44cb-bf6b-eed6400b361f.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.zh-cn.xlf"
$ws.Range("E2").Value = "2016-03-24 03:14:36"
$ws.Range("F2").Value = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.md"
$ws.Range("G2").Value = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.zh-cn.xlf"
$ws.Range("H2").Value = "2016-03-24 03:14:59"
$ws.Range("J2").Value = "Include"
$ws.Range("A3").Value = "ffffffdf71933c-d945-47f1-b0b7-c48448e4a22c.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-24 03:14:36"
$ws.Range("F3").Value = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.md"
$ws.Range("G3").Value = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.zh-cn.xlf"
$ws.Range("H3").Value = "2016-03-24 03:14:59"
$ws.Range("J3").Value = "Include"
$ws.Range("A4").Value = "9ac5998c-9398-47e2-b777-121c530ec423.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "9ac5998c-9398-47e2-b777-121c530ec423.c2f1658e80855fe2ad137393932fd933c9d75be2.zh-cn.xlf"
$ws.Range("E4").Value = "2016-03-24 03:19:22"
$ws.Range("F4").Value = "9ac5998c-9398-47e2-b777-121c530ec423.md"
$ws.Range("G4").Value = "9ac5998c-9398-47e2-b777-121c530ec423.c2f1658e80855fe2ad137393932fd933c9d75be2.zh-cn.xlf"
$ws.Range("H4").Value = "2016-03-24 03:18:32"
$ws.Range("J4").Value = "Include"

# Update hyperlink display text in place (keep existing target URL/rId)
foreach ($h in $ws.Hyperlinks) {
  $addr = $h.Range.Address()
  if ($addr -eq '$A$2') { $h.TextToDisplay = "ffff96a6b1e0-7098-44cb-bf6b-eed6400b361f.md" }
  if ($addr -eq '$D$2') { $h.TextToDisplay = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.zh-cn.xlf" }
  if ($addr -eq '$F$2') { $h.TextToDisplay = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.md" }
  if ($addr -eq '$G$2') { $h.TextToDisplay = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.zh-cn.xlf" }
  if ($addr -eq '$A$3') { $h.TextToDisplay = "ffffffdf71933c-d945-47f1-b0b7-c48448e4a22c.md" }
  if ($addr -eq '$D$3') { $h.TextToDisplay = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.zh-cn.xlf" }
  if ($addr -eq '$F$3') { $h.TextToDisplay = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.md" }
  if ($addr -eq '$G$3') { $h.TextToDisplay = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.zh-cn.xlf" }
  if ($addr -eq '$A$4') { $h.TextToDisplay = "9ac5998c-9398-47e2-b777-121c530ec423.md" }
  if ($addr -eq '$D$4') { $h.TextToDisplay = "9ac5998c-9398-47e2-b777-121c530ec423.c2f1658e80855fe2ad137393932fd933c9d75be2.zh-cn.xlf" }
  if ($addr -eq '$F$4') { $h.TextToDisplay = "9ac5998c-9398-47e2-b777-121c530ec423.md" }
  if ($addr -eq '$G$4') { $h.TextToDisplay = "9ac5998c-9398-47e2-b777-121c530ec423.c2f1658e80855fe2ad137393932fd933c9d75be2.zh-cn.xlf" }
}

# ---- Sheet: de-de ----
$ws = $wb.Worksheets.Item("de-de")

# Update cell values (data rows 2-4; header row unchanged)
$ws.Range("A2").Value = "ffff96a6b1e0-7098-44cb-bf6b-eed6400b361f.md"
$ws.Range("B2").Value = ".md"
$ws.Range("C2").Value = "Handed back: in sync with en-US"
$ws.Range("D2").Value = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.de-de.xlf"
$ws.Range("E2").Value = "2016-03-24 03:14:40"
$ws.Range("F2").Value = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.md"
$ws.Range("G2").Value = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.de-de.xlf"
$ws.Range("H2").Value = "2016-03-24 03:15:08"
$ws.Range("J2").Value = "Include"
$ws.Range("A3").Value = "ffffffdf71933c-d945-47f1-b0b7-c48448e4a22c.md"
$ws.Range("B3").Value = ".md"
$ws.Range("C3").Value = "Handed back: in sync with en-US"
$ws.Range("D3").Value = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.de-de.xlf"
$ws.Range("E3").Value = "2016-03-24 03:14:40"
$ws.Range("F3").Value = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.md"
$ws.Range("G3").Value = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.de-de.xlf"
$ws.Range("H3").Value = "2016-03-24 03:15:08"
$ws.Range("J3").Value = "Include"
$ws.Range("A4").Value = "9ac5998c-9398-47e2-b777-121c530ec423.md"
$ws.Range("B4").Value = ".md"
$ws.Range("C4").Value = "Ready for handoff"
$ws.Range("D4").Value = "9ac5998c-9398-47e2-b777-121c530ec423.c2f1658e80855fe2ad137393932fd933c9d75be2.de-de.xlf"
$ws.Range("E4").Value = "2016-03-24 03:19:26"
$ws.Range("F4").Value = "9ac5998c-9398-47e2-b777-121c530ec423.md"
$ws.Range("G4").Value = "9ac5998c-9398-47e2-b777-121c530ec423.c2f1658e80855fe2ad137393932fd933c9d75be2.de-de.xlf"
$ws.Range("H4").Value = "2016-03-24 03:18:39"
$ws.Range("J4").Value = "Include"

# Update hyperlink display text in place (keep existing target URL/rId)
foreach ($h in $ws.Hyperlinks) {
  $addr = $h.Range.Address()
  if ($addr -eq '$A$2') { $h.TextToDisplay = "ffff96a6b1e0-7098-44cb-bf6b-eed6400b361f.md" }
  if ($addr -eq '$D$2') { $h.TextToDisplay = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.de-de.xlf" }
  if ($addr -eq '$F$2') { $h.TextToDisplay = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.md" }
  if ($addr -eq '$G$2') { $h.TextToDisplay = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.de-de.xlf" }
  if ($addr -eq '$A$3') { $h.TextToDisplay = "ffffffdf71933c-d945-47f1-b0b7-c48448e4a22c.md" }
  if ($addr -eq '$D$3') { $h.TextToDisplay = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.de-de.xlf" }
  if ($addr -eq '$F$3') { $h.TextToDisplay = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.md" }
  if ($addr -eq '$G$3') { $h.TextToDisplay = "0e8cf8ef-95db-4b11-be8d-95caf4d79894.071e4d04458edb58284939ebb9be732bdf773d2e.de-de.xlf" }
  if ($addr -eq '$A$4') { $h.TextToDisplay = "9ac5998c-9398-47e2-b777-121c530ec423.md" }
  if ($addr -eq '$D$4') { $h.TextToDisplay = "9ac5998c-9398-47e2-b777-121c530ec423.c2f1658e80855fe2ad137393932fd933c9d75be2.de-de.xlf" }
  if ($addr -eq '$F$4') { $h.TextToDisplay = "9ac5998c-9398-47e2-b777-121c530ec423.md" }
  if ($addr -eq '$G$4') { $h.TextToDisplay = "9ac5998c-9398-47e2-b777-121c530ec423.c2f1658e80855fe2ad137393932fd933c9d75be2.de-de.xlf" }
}
